$wb = $excel.ActiveWorkbook

# Update the "Estado" (Status) column on the Backlog sheet.
$backlog = $wb.Worksheets.Item("Backlog")
$backlog.Range("D4").Value = "Realizado"
$backlog.Range("D5").Value = "Realizado"
$backlog.Range("D7").Value = "En proceso"
$backlog.Range("D8").Value = "Realizado"

# Make Backlog the active sheet, with D9 selected.
$backlog.Activate()
$backlog.Range("D9").Select()

# The Instructivo sheet should no longer be the tab that is selected.
